# Auto-generated edit script applying the Spriggan_Profits data refresh
# Updates numeric cached values in columns H-N for specific rows across
# the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2046.6
$ws.Range("I80").Value = 3591.8
$ws.Range("J80").Value = 501.4
$ws.Range("K80").Value = 10775.4
$ws.Range("L80").Value = 1504.2
$ws.Range("M80").Value = -9777.400000000001
$ws.Range("N80").Value = -3500.2
$ws.Range("H83").Value = 2046.6
$ws.Range("I83").Value = 3591.8
$ws.Range("J83").Value = 501.4
$ws.Range("K83").Value = 32326.2
$ws.Range("L83").Value = 4512.599999999999
$ws.Range("M83").Value = -27334.2
$ws.Range("N83").Value = -14496.6
$ws.Range("H92").Value = 764.86365
$ws.Range("I92").Value = 722.3684
$ws.Range("K92").Value = 722.3684
$ws.Range("M92").Value = 525.6316
$ws.Range("H98").Value = 1147.8857
$ws.Range("I98").Value = 887.5294
$ws.Range("K98").Value = 887.5294
$ws.Range("M98").Value = 610.4706
$ws.Range("H122").Value = 1147.8857
$ws.Range("I122").Value = 887.5294
$ws.Range("K122").Value = 2662.5882
$ws.Range("M122").Value = -212.5882000000001
$ws.Range("H127").Value = 942.3
$ws.Range("I127").Value = 942.3
$ws.Range("K127").Value = 2826.9
$ws.Range("M127").Value = 2133.1
$ws.Range("H129").Value = 9752.799999999999
$ws.Range("I129").Value = 1366.375
$ws.Range("K129").Value = 4099.125
$ws.Range("M129").Value = 900.875
$ws.Range("H131").Value = 795.5833
$ws.Range("I131").Value = 795.5833
$ws.Range("K131").Value = 2386.7499
$ws.Range("M131").Value = 2653.2501
$ws.Range("H132").Value = 3283.05
$ws.Range("I132").Value = 3401.0667
$ws.Range("J132").Value = 2929
$ws.Range("K132").Value = 10203.2001
$ws.Range("L132").Value = 8787
$ws.Range("M132").Value = -7673.2001
$ws.Range("N132").Value = -13847
$ws.Range("H137").Value = 2182.24
$ws.Range("I137").Value = 1394.7059
$ws.Range("J137").Value = 3855.75
$ws.Range("K137").Value = 4184.1177
$ws.Range("L137").Value = 11567.25
$ws.Range("M137").Value = -1634.1177
$ws.Range("N137").Value = -16667.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1936.1111
$ws.Range("I32").Value = 1954.0698
$ws.Range("J32").Value = 1550
$ws.Range("K32").Value = 1954.0698
$ws.Range("L32").Value = 1550
$ws.Range("M32").Value = -1667.0698
$ws.Range("N32").Value = -2124
$ws.Range("H97").Value = 134.4
$ws.Range("I97").Value = 138
$ws.Range("J97").Value = 120
$ws.Range("K97").Value = 138
$ws.Range("L97").Value = 120
$ws.Range("M97").Value = 358
$ws.Range("N97").Value = -1112

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 50000
$ws.Range("J61").Value = 50000
$ws.Range("L61").Value = 50000
$ws.Range("N61").Value = -50626
$ws.Range("H86").Value = 3274.516
$ws.Range("I86").Value = 3500.7144
$ws.Range("K86").Value = 3500.7144
$ws.Range("M86").Value = -2377.7144
$ws.Range("H89").Value = 3274.516
$ws.Range("I89").Value = 3500.7144
$ws.Range("K89").Value = 17503.572
$ws.Range("M89").Value = -11887.572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2232.2222
$ws.Range("I16").Value = 2159.6365
$ws.Range("J16").Value = 2346.2856
$ws.Range("K16").Value = 2159.6365
$ws.Range("L16").Value = 2346.2856
$ws.Range("M16").Value = -1872.6365
$ws.Range("N16").Value = -2920.2856
$ws.Range("H22").Value = 36833
$ws.Range("J22").Value = 5250
$ws.Range("L22").Value = 5250
$ws.Range("N22").Value = -5950
$ws.Range("H31").Value = 6297.7144
$ws.Range("I31").Value = 3651.6365
$ws.Range("J31").Value = 16000
$ws.Range("K31").Value = 3651.6365
$ws.Range("L31").Value = 16000
$ws.Range("M31").Value = -3356.6365
$ws.Range("N31").Value = -16590
$ws.Range("H34").Value = 6297.7144
$ws.Range("I34").Value = 3651.6365
$ws.Range("J34").Value = 16000
$ws.Range("K34").Value = 3651.6365
$ws.Range("L34").Value = 16000
$ws.Range("M34").Value = -3449.6365
$ws.Range("N34").Value = -16404
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = $null
$ws.Range("N47").Value = $null
$ws.Range("H54").Value = 33599.6
$ws.Range("J54").Value = 33599.6
$ws.Range("L54").Value = 33599.6
$ws.Range("N54").Value = -34915.6
$ws.Range("H58").Value = 23819200
$ws.Range("I58").Value = 23819200
$ws.Range("K58").Value = 23819200
$ws.Range("M58").Value = -23818997
$ws.Range("H94").Value = 1695.4
$ws.Range("J94").Value = 1681.2222
$ws.Range("L94").Value = 1681.2222
$ws.Range("N94").Value = -2583.2222
$ws.Range("H97").Value = 29665
$ws.Range("J97").Value = 31598.2
$ws.Range("L97").Value = 31598.2
$ws.Range("N97").Value = -33580.2
$ws.Range("H113").Value = 2232.2222
$ws.Range("I113").Value = 2159.6365
$ws.Range("J113").Value = 2346.2856
$ws.Range("K113").Value = 2159.6365
$ws.Range("L113").Value = 2346.2856
$ws.Range("M113").Value = 10.36349999999993
$ws.Range("N113").Value = -6686.2856
$ws.Range("H134").Value = 27780818
$ws.Range("I134").Value = 27780818
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 83342454
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -83339919
$ws.Range("N134").Value = $null
$ws.Range("H136").Value = 23819200
$ws.Range("I136").Value = 23819200
$ws.Range("K136").Value = 71457600
$ws.Range("M136").Value = -71455050

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 3995
$ws.Range("I57").Value = 3995
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 11985
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -11426
$ws.Range("N57").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 5600
$ws.Range("I46").Value = 4500
$ws.Range("K46").Value = 4500
$ws.Range("M46").Value = -4344
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").Value = $null
$ws.Range("H70").Value = 7283.625
$ws.Range("I70").Value = 7249.8335
$ws.Range("J70").Value = 7385
$ws.Range("K70").Value = 7249.8335
$ws.Range("L70").Value = 7385
$ws.Range("M70").Value = -6979.8335
$ws.Range("N70").Value = -7925
$ws.Range("H73").Value = 7283.625
$ws.Range("I73").Value = 7249.8335
$ws.Range("J73").Value = 7385
$ws.Range("K73").Value = 7249.8335
$ws.Range("L73").Value = 7385
$ws.Range("M73").Value = -6313.8335
$ws.Range("N73").Value = -9257
$ws.Range("H80").Value = 3664.1667
$ws.Range("I80").Value = 3664.1667
$ws.Range("K80").Value = 3664.1667
$ws.Range("M80").Value = -2666.1667
$ws.Range("H83").Value = 3664.1667
$ws.Range("I83").Value = 3664.1667
$ws.Range("K83").Value = 18320.8335
$ws.Range("M83").Value = -13328.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1462.1154
$ws.Range("J136").Value = 1983.2
$ws.Range("L136").Value = 5949.6
$ws.Range("N136").Value = -11049.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 25000
$ws.Range("J51").Value = 25000
$ws.Range("L51").Value = 25000
$ws.Range("N51").Value = -26020
$ws.Range("H112").Value = 39166.332
$ws.Range("J112").Value = 39166.332
$ws.Range("L112").Value = 39166.332
$ws.Range("N112").Value = -42120.332
